$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder / tweak maa:// list strings (star-prefix counts shuffled) ---
$ws.Range("L2").Value = "*maa://24633, maa://39402, *maa://30515, *maa://34787, ***maa://29083"
$ws.Range("T3").Value = "maa://24617, maa://45854, **maa://20790"
$ws.Range("D12").Value = "maa://30766, maa://36678"
$ws.Range("H12").Value = "maa://21867, **maa://45826"
$ws.Range("L18").Value = "maa://22466, **maa://22732"

# --- Update refresh timestamp (row 8, column A) ---
$ws.Range("A8").Value = "更新日期：2025.05.02 13:20:57"

# --- Row 31: 小满 gained a recorded drop ---
# Leading "'" keeps the digit-only text as text (matches source inlineStr "1")
# instead of Excel's normal auto-coercion to a number.
$ws.Range("AA31").Value = "'1"
$ws.Range("AB31").Value = "**maa://51420"

# --- Row 34: 芳汀 / 诺威尔 now have recorded (empty/"None") results instead of "-" ---
$ws.Range("G34").Value = "'0"
$ws.Range("H34").Value = "None"
$ws.Range("W34").Value = "'0"
$ws.Range("X34").Value = "None"

# --- Row 37: new operator block added at B:E (CONFESS-47) ---
$ws.Range("B37").Value = "CONFESS-47"
$ws.Range("C37").Value = "-"
$ws.Range("D37").Value = "-"

# --- Row 43: new operator block added at J:M (信仰搅拌机) ---
$ws.Range("J43").Value = "信仰搅拌机"
$ws.Range("K43").Value = "-"
$ws.Range("L43").Value = "-"

# --- Row 45: new operator block added at AD:AG (新约能天使) ---
$ws.Range("AD45").Value = "新约能天使"
$ws.Range("AE45").Value = "-"
$ws.Range("AF45").Value = "-"

# --- Row 54: new operator block added at N:Q (蕾缪安) ---
$ws.Range("N54").Value = "蕾缪安"
$ws.Range("O54").Value = "-"
$ws.Range("P54").Value = "-"

# --- New row 76: new operator block added at F:I (聆音) ---
$ws.Range("F76").Value = "聆音"
$ws.Range("G76").Value = "-"
$ws.Range("H76").Value = "-"
